# Update the 8 site/location code values in column B: append the "_co"
# country suffix (e.g. "Maaji_Retiro" -> "Maaji_Retiro_co") and drop the
# (redundant/default) explicit cell style that had been applied to them.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2:B9").ClearFormats()

$ws.Range("B2").Value = "Maaji_Retiro_co"
$ws.Range("B3").Value = "Maaji_Barranquilla_co"
$ws.Range("B4").Value = "Maaji_Cartagena_co"
$ws.Range("B5").Value = "Maaji_Tesoro_co"
$ws.Range("B6").Value = "Maaji_Cali_co"
$ws.Range("B7").Value = "Maaji_Bocagrande_co"
$ws.Range("B8").Value = "Maaji_Santafe_co"
$ws.Range("B9").Value = "Maaji_Bucaramanga_co"

# The workbook shipped with a full set of unused built-in cell styles
# (Good/Bad/Neutral, accent colors, etc.) inherited from the default
# template; none of them were actually applied to any cell (only the
# implicit "Normal" style was in use). Remove them so the style table
# collapses back down to just "Normal".
$styleNames = @()
foreach ($style in $wb.Styles) {
    $styleNames += $style.Name
}
foreach ($name in $styleNames) {
    if ($name -ne "Normal") {
        $wb.Styles.Item($name).Delete()
    }
}

# Scroll/selection state: the sheet was left with V2:V9 selected.
$ws.Range("V2:V9").Select()
